$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value for the completed 2.13 Switch statement section
$ws.Range("A18").Value = 2.13

# Update the selection to match the new active cell after the edit
$ws.Range("A18").Select()
